# Updated map March 2024
# Adds three new beer rows (87-89) to the beer-names map and widens
# column A so the new (longer) names display nicely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column A (new "name key" column) -------------------------------
# Target OOXML width is 27.81; Excel/this COM layer quantizes ColumnWidth to
# whole-pixel steps, so 27 is the closest input that lands on the nearest
# achievable stored width.
$ws.Columns.Item(1).ColumnWidth = 27

# --- Append the three new rows ---------------------------------------------
$ws.Cells.Item(87, 1).Value2 = "WHITE MOUNTAIN WHITE ALE"
$ws.Cells.Item(87, 2).Value2 = "White Mountain Ale"
$ws.Cells.Item(87, 3).Value2 = "Other"

$ws.Cells.Item(88, 1).Value2 = "SQUIRREL FIGHTS NUT BROWN"
$ws.Cells.Item(88, 2).Value2 = "Squirrel Fights Nut Brown Ale"
$ws.Cells.Item(88, 3).Value2 = "Other"

$ws.Cells.Item(89, 1).Value2 = "OKTOBERFEST MARZEN"
$ws.Cells.Item(89, 2).Value2 = "Oktoberfest"
$ws.Cells.Item(89, 3).Value2 = "Other"

# --- Restore the view state (selection moved to D75 while editing) --------
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("D75").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 44
$excel.ActiveWindow.ScrollColumn = 1
